# Bulk Upload Final changes
# Replace the 30 proctor records in rows 2-31 (columns A-D) with a new
# batch of generated Proctor/Automation accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Proctor66783", "Automation66783", "proctorautomation66783@gmail.com", "66783"),
    @("Proctor20305", "Automation20305", "proctorautomation20305@gmail.com", "20305"),
    @("Proctor35557", "Automation35557", "proctorautomation35557@gmail.com", "35557"),
    @("Proctor69732", "Automation69732", "proctorautomation69732@gmail.com", "69732"),
    @("Proctor80170", "Automation80170", "proctorautomation80170@gmail.com", "80170"),
    @("Proctor25107", "Automation25107", "proctorautomation25107@gmail.com", "25107"),
    @("Proctor56119", "Automation56119", "proctorautomation56119@gmail.com", "56119"),
    @("Proctor55471", "Automation55471", "proctorautomation55471@gmail.com", "55471"),
    @("Proctor19851", "Automation19851", "proctorautomation19851@gmail.com", "19851"),
    @("Proctor51868", "Automation51868", "proctorautomation51868@gmail.com", "51868"),
    @("Proctor33137", "Automation33137", "proctorautomation33137@gmail.com", "33137"),
    @("Proctor57154", "Automation57154", "proctorautomation57154@gmail.com", "57154"),
    @("Proctor62732", "Automation62732", "proctorautomation62732@gmail.com", "62732"),
    @("Proctor76947", "Automation76947", "proctorautomation76947@gmail.com", "76947"),
    @("Proctor19887", "Automation19887", "proctorautomation19887@gmail.com", "19887"),
    @("Proctor87654", "Automation87654", "proctorautomation87654@gmail.com", "87654"),
    @("Proctor70339", "Automation70339", "proctorautomation70339@gmail.com", "70339"),
    @("Proctor06015", "Automation06015", "proctorautomation06015@gmail.com", "06015"),
    @("Proctor79412", "Automation79412", "proctorautomation79412@gmail.com", "79412"),
    @("Proctor89009", "Automation89009", "proctorautomation89009@gmail.com", "89009"),
    @("Proctor52298", "Automation52298", "proctorautomation52298@gmail.com", "52298"),
    @("Proctor70595", "Automation70595", "proctorautomation70595@gmail.com", "70595"),
    @("Proctor60576", "Automation60576", "proctorautomation60576@gmail.com", "60576"),
    @("Proctor47712", "Automation47712", "proctorautomation47712@gmail.com", "47712"),
    @("Proctor76108", "Automation76108", "proctorautomation76108@gmail.com", "76108"),
    @("Proctor06476", "Automation06476", "proctorautomation06476@gmail.com", "06476"),
    @("Proctor97659", "Automation97659", "proctorautomation97659@gmail.com", "97659"),
    @("Proctor74639", "Automation74639", "proctorautomation74639@gmail.com", "74639"),
    @("Proctor30431", "Automation30431", "proctorautomation30431@gmail.com", "30431"),
    @("Proctor03337", "Automation03337", "proctorautomation03337@gmail.com", "03337")
)

$firstRow = 2
$lastRow = $firstRow + $data.Count - 1

# Column D holds numeric-looking IDs that may include leading zeros
# (e.g. "06015"). Temporarily mark the column as Text so Excel stores
# the literal digit string instead of coercing it to a number, then
# restore the default formatting once the values are in place.
$idRange = $ws.Range("D" + $firstRow + ":D" + $lastRow)
$idRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$idRange.ClearFormats()

